# 4.0.3 model and data
#
# Splits the single "trans/BVTQaZ/BVTQaZ.csv" and "trans/VTQaZ/VTQaZ.csv"
# rows on the "Boolean" sheet into six mode-specific files each (LDVs, HDVs,
# aircraft, rail, ships, motorbikes), and updates the active-sheet /
# selection state left behind by the author's last save.

$wb = $excel.ActiveWorkbook

$wsAbout     = $wb.Worksheets.Item("About")
$wsInteger   = $wb.Worksheets.Item("Integer")
$wsBoolean   = $wb.Worksheets.Item("Boolean")
$wsSubscript = $wb.Worksheets.Item("Subscript")

# ---------------------------------------------------------------------
# Boolean sheet: expand the two combined transport CSVs into six rows
# apiece, inserted in place of the single summary row.
# ---------------------------------------------------------------------

$bvtqazSplit = @(
  "trans/BVTQaZ/BVTQaZ-LDVs.csv",
  "trans/BVTQaZ/BVTQaZ-HDVs.csv",
  "trans/BVTQaZ/BVTQaZ-aircraft.csv",
  "trans/BVTQaZ/BVTQaZ-rail.csv",
  "trans/BVTQaZ/BVTQaZ-ships.csv",
  "trans/BVTQaZ/BVTQaZ-motorbikes.csv"
)

$vtqazSplit = @(
  "trans/VTQaZ/VTQaZ-LDVs.csv",
  "trans/VTQaZ/VTQaZ-HDVs.csv",
  "trans/VTQaZ/VTQaZ-aircraft.csv",
  "trans/VTQaZ/VTQaZ-rail.csv",
  "trans/VTQaZ/VTQaZ-ships.csv",
  "trans/VTQaZ/VTQaZ-motorbikes.csv"
)

# "trans/BVTQaZ/BVTQaZ.csv" currently sits on row 17. Make room for five
# more rows (6 total), then write the six split filenames.
$wsBoolean.Rows("17:21").Insert()
for ($i = 0; $i -lt $bvtqazSplit.Length; $i++) {
  $wsBoolean.Range("A" + (17 + $i)).Value = $bvtqazSplit[$i]
}

# "trans/VTQaZ/VTQaZ.csv" used to be row 21; after the six-row insert above
# it is now row 26. Make room for five more rows, then write the split.
$wsBoolean.Rows("26:30").Insert()
for ($i = 0; $i -lt $vtqazSplit.Length; $i++) {
  $wsBoolean.Range("A" + (26 + $i)).Value = $vtqazSplit[$i]
}

# ---------------------------------------------------------------------
# View/selection state: the author's last save left "About" as the active
# tab, with a lingering selection on "Integer" (A13) and "Boolean"
# (scrolled to row 10, selection on A32).
# ---------------------------------------------------------------------

$wsInteger.Range("A13").Select()
$wsBoolean.Range("A32").Select()
$wsBoolean.Application.ActiveWindow.ScrollRow = 10

$wsAbout.Activate()
$wsAbout.Range("A1").Select()
